{"js": "// Add two new paragraphs at the end of the document body:\n//   1) \"2022\u5e746\u670810\u65e5\u661f\u671f\u4e94\" (split across 4 runs, matching how Word\n//      would naturally segment mixed Latin-digit / East-Asian-hinted text)\n//   2) \"\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u65e2\u7b80\u5355\u53c8\u4fbf\u6377\u3002\"\n//\n// We build the exact run structure with insertOoxml (FlatOpc) so that the\n// runs keep their distinct <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n// markings instead of being silently coalesced into a single run (which is\n// what happens if the same text is inserted via insertParagraph/insertText,\n// since those calls merge adjacent runs that end up with identical\n// formatting).\n\nconst body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\nconst insertionPoint = lastParagraph.getRange(Word.RangeLocation.after);\n\nconst flatOpcXml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>2022</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u5e746\u67081</w:t></w:r>\n            <w:r><w:t>0</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u65e5\u661f\u671f\u4e94</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u65e2\u7b80\u5355\u53c8\u4fbf\u6377\u3002</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add two new paragraphs at the end of the document body:\n#   1) \"2022\u5e746\u670810\u65e5\u661f\u671f\u4e94\" (split across 4 runs, matching how Word\n#      naturally segments mixed Latin-digit / East-Asian-hinted text)\n#   2) \"\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u65e2\u7b80\u5355\u53c8\u4fbf\u6377\u3002\"\n#\n# Range.InsertXML REPLACES the target range's contents, so we first create a\n# fresh, empty trailing paragraph to host the inserted XML (keeping the\n# original last paragraph intact), then clean up the extra empty paragraph\n# mark that InsertXML leaves behind afterwards.\n\n$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n$tailRange = $lastPara.Range\n$tailRange.Collapse(0)          # wdCollapseEnd\n[void]$tailRange.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$insertionRange = $newPara.Range\n$insertionRange.Collapse(0)     # wdCollapseEnd\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>2022</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u5e746\u67081</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u65e5\u661f\u671f\u4e94</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u65e2\u7b80\u5355\u53c8\u4fbf\u6377\u3002</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n[void]$insertionRange.InsertXML($xml)\n\n# InsertXML leaves the paragraph mark that was there before it behind as a\n# trailing empty paragraph. Remove it by deleting the paragraph mark of the\n# paragraph that now holds our newly inserted text, which merges the\n# (empty) final paragraph back into it.\n$count = $d.Paragraphs.Count\n$markRange = $d.Paragraphs.Item($count - 1).Range\n$markRange.Collapse(0)          # wdCollapseEnd\n[void]$markRange.MoveEnd(1, 1)  # wdCharacter: extend over the paragraph mark\n[void]$markRange.Delete()\n"}
